# Fixed issue #13 Permitir que en los ficheros de metadatos dos columnas se
# puedan relacionar para crear SKOS jerarquicos.
#
# This inserts a new metadata row right below the header row. The new row
# holds a short "slug" identifier for every column (used to relate two
# columns together, e.g. "provincia-nombre" / "provincia-codigo"), pushing
# the previous measure/dimension metadata rows down by one. The now
# orphaned last row (which only held the stray "mapping-ano.xlsx" value)
# is dropped in the process.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 2-5 down to 3-6 to make room for the new slug row.
$ws.Rows.Item(2).Insert()

# New row 2: short slug identifiers for each column.
$ws.Range("A2").Value = "poblacion"
$ws.Range("B2").Value = "ccaa-nombre"
$ws.Range("C2").Value = "abstencion"
$ws.Range("D2").Value = "comarca-codigo"
$ws.Range("E2").Value = "provincia-nombre"
$ws.Range("F2").Value = "municipio-nombre"
$ws.Range("G2").Value = "votos-blancos"
$ws.Range("H2").Value = "votos-nulos"
$ws.Range("I2").Value = "comarca-nombre"
$ws.Range("J2").Value = "participacion"
$ws.Range("K2").Value = "ccaa-codigo"
$ws.Range("L2").Value = "participacion"
$ws.Range("M2").Value = "censo-electoral"
$ws.Range("N2").Value = "votos-a-candidaturas"
$ws.Range("O2").Value = "provincia-codigo"
$ws.Range("P2").Value = "municipio-codigo"
$ws.Range("Q2").Value = "ano"

# The row that used to be row 5 (now row 6) only carried a stray
# "mapping-ano.xlsx" leftover value in Q; it is no longer needed now
# that every column has its own slug identifier, so remove it entirely.
$ws.Rows.Item(6).Delete()
